$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top; this shifts all existing data down by 2 rows,
# turning the old "A1:C79" data range into "A3:C81".
$ws.Rows("1:2").Insert()

# Row 1: new header labels
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"

# Row 2: column names
$ws.Range("A2").Value = "municipio"
$ws.Range("B2").Value = "CASOS"
$ws.Range("C2").Value = "óbitos"

# Style row 1: bold font, thin box border around every cell, centered horizontally
# and aligned to the top vertically. Build the style once on A1 and fan it out to
# B1:C1 with a format-only paste so the workbook ends up with a single extra style
# record instead of one per incremental property write.
$a1 = $ws.Range("A1")
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1

$a1.Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
